$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.878.32"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "'1.806.48"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").Value = "'231.46"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "'38.91"
$ws.Range("E8").Value = "  -8.72%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").Value = "'0.0678"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("D11").Value = "'0.0993"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'2.066.46"
$ws.Range("D13").Value = "'1.798.23"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").Value = "'10.88"
$ws.Range("E15").Value = "  -6.11%  "
$ws.Range("E16").Value = "  -5.33%  "
$ws.Range("D17").Value = "'34.794.33"
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("D18").Value = "'69.03"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").Value = "0.0₃0779"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").Value = "'237.96"
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("D21").Value = "'11.76"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'171.58"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "'7.75"
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").Value = "'17.23"
$ws.Range("E27").Value = "  -4.24%  "
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("D29").Value = "'1.55"
$ws.Range("E29").Value = "  +6.49%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "'3.97"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("D34").Value = "'1.76"
$ws.Range("E34").Value = "  -8.19%  "
$ws.Range("D35").Value = "'1.14"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("D36").Value = "'0.680"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "'90.48"
$ws.Range("E37").Value = "  -8.89%  "
$ws.Range("D38").Value = "'1.31"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'1.304.65"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0191"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("E42").Value = "  -5.22%  "
$ws.Range("D43").Value = "'0.957"
$ws.Range("E43").Value = "  -7.03%  "
$ws.Range("D44").Value = "'2.20"
$ws.Range("E44").Value = "  -13.22%  "
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").Value = "'6.13"
$ws.Range("D47").Value = "'0.0509"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").Value = "'1.986.40"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0672"
$ws.Range("E49").Value = "  +7.34%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.01"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "'98.80"
$ws.Range("E51").Value = "  +2.11%  "
